$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "26.027.90"
$ws.Range("E2").Value = "  -2.03%  "
Set-TextValue "D3" "1.668.15"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  -0.15%  "
Set-TextValue "D5" "216.73"
$ws.Range("E5").Value = "  -1.28%  "
Set-TextValue "D6" "0.5095"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  -0.18%  "
Set-TextValue "D8" "0.2654"
$ws.Range("E8").Value = "  -0.10%  "
Set-TextValue "D9" "0.06389"
$ws.Range("E9").Value = "  +1.35%  "
Set-TextValue "D10" "21.79"
$ws.Range("E10").Value = "  -0.93%  "
Set-TextValue "D11" "0.07455"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D12" "4.511"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.652.56"
$ws.Range("E13").Value = "  -2.27%  "
Set-TextValue "D14" "0.5824"
$ws.Range("E14").Value = "  +0.88%  "
Set-TextValue "D15" "0.000008550"
$ws.Range("E15").Value = "  +0.65%  "
Set-TextValue "D16" "64.33"
$ws.Range("E16").Value = "  -1.30%  "
Set-TextValue "D17" "26.129.67"
$ws.Range("E17").Value = "  -1.74%  "
Set-TextValue "D18" "4.937"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("E20").Value = "  -1.37%  "
Set-TextValue "D21" "191.40"
$ws.Range("E21").Value = "  +2.78%  "
Set-TextValue "D22" "6.197"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("E23").Value = "  -0.13%  "
Set-TextValue "D24" "144.70"
$ws.Range("E24").Value = "  -0.04%  "
Set-TextValue "D25" "7.609"
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("E26").Value = "  +2.59%  "
$ws.Range("E27").Value = "  -0.72%  "
Set-TextValue "D28" "0.06586"
$ws.Range("E28").Value = "  +15.28%  "
Set-TextValue "D29" "1.339"
$ws.Range("E29").Value = "  +0.25%  "
Set-TextValue "D30" "1.316"
$ws.Range("E30").Value = "  -1.23%  "
Set-TextValue "D31" "3.542"
$ws.Range("E31").Value = "  +0.75%  "
Set-TextValue "D32" "3.512"
$ws.Range("E32").Value = "  +0.38%  "
Set-TextValue "D33" "1.654"
$ws.Range("E33").Value = "  +0.22%  "
Set-TextValue "D34" "1.018"
$ws.Range("E34").Value = "  -0.24%  "
Set-TextValue "D35" "0.6121"
$ws.Range("E35").Value = "  +2.36%  "
Set-TextValue "D36" "2.369"
$ws.Range("E36").Value = "  +0.15%  "
Set-TextValue "D37" "2.694"
$ws.Range("E37").Value = "  +0.64%  "
Set-TextValue "D38" "6.294"
$ws.Range("E38").Value = "  +8.01%  "
Set-TextValue "D39" "1.092.22"
$ws.Range("E39").Value = "  +0.31%  "
Set-TextValue "D40" "0.01600"
$ws.Range("E40").Value = "  -1.22%  "
Set-TextValue "D41" "0.8732"
$ws.Range("E41").Value = "  +1.35%  "
Set-TextValue "D42" "1.009"
$ws.Range("E42").Value = "  +0.35%  "
Set-TextValue "D43" "101.10"
$ws.Range("E43").Value = "  +1.45%  "
Set-TextValue "D44" "1.816.61"
$ws.Range("E44").Value = "  -1.77%  "
Set-TextValue "D45" "0.00000000114"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  +0.38%  "
Set-TextValue "D48" "8.076"
$ws.Range("E48").Value = "  -0.19%  "
Set-TextValue "D49" "0.05231"
$ws.Range("E49").Value = "  -0.06%  "
Set-TextValue "D50" "0.4288"
$ws.Range("E50").Value = "  -0.81%  "
Set-TextValue "D51" "6.027"
$ws.Range("E51").Value = "  +4.22%  "
